$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Restricciones_del_lider (sheet index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "0.8 - x"
Set-TextValue $ws "B2" "-1.8"
Set-TextValue $ws "D2" "0.74"
$ws.Range("A3").Value = "-0.8 + x"
Set-TextValue $ws "B3" "-0.19999999999999996"
Set-TextValue $ws "D3" "0.96"

# Restricciones_del_follower (sheet index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "-4.473856209150328 + 2.4183006535947715y"
Set-TextValue $ws "B2" "3.473856209150328"
Set-TextValue $ws "D2" "0.76"
Set-TextValue $ws "E2" "3.0"
Set-TextValue $ws "F2" "3.7"
$ws.Range("A3").Value = "0.9805000000000001 - 0.53y"
Set-TextValue $ws "B3" "-1.9805000000000001"
Set-TextValue $ws "D3" "0.2"
Set-TextValue $ws "E3" "0"
Set-TextValue $ws "F3" "9.1"

# Punto_modificado (sheet index 4)
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws "A2" "0.8"
Set-TextValue $ws "B2" "1.85"

# Vector_bf (sheet index 5)
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws "A2" "-11.435158496732027"

# Vector_BF (sheet index 6)
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws "A2" "4.591"
Set-TextValue $ws "A3" "-74.02545196078431"

# Vector_Alpha (sheet index 7) - plain numeric cell, not text
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 1.53
